$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that was updated
# from 45170 (2023-09-01) to 45174 (2023-09-05) for rows 2 through 15.
$ws.Range("C2:C15").Value = 45174
